$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: 1752. Check if Array Is Sorted and Rotated (Java)
$ws.Cells.Item(21, 1).Value = 1752
$ws.Cells.Item(21, 2).Value = "Check if Array Is Sorted and Rotated"
$ws.Cells.Item(21, 2).Style = "Normal"
$ws.Cells.Item(21, 3).Value = "Java"
$ws.Cells.Item(21, 4).Value = 45000
$ws.Cells.Item(21, 4).NumberFormat = "d-mmm-yy"

# Row 22: 26. Remove Duplicates from Sorted Array (Java)
$ws.Cells.Item(22, 1).Value = 26
$ws.Cells.Item(22, 2).Value = "Remove Duplicates from Sorted Array"
$ws.Cells.Item(22, 2).Style = "Normal"
$ws.Cells.Item(22, 3).Value = "Java"
$ws.Cells.Item(22, 4).Value = 45000
$ws.Cells.Item(22, 4).NumberFormat = "d-mmm-yy"

# Update the view: clear the scrolled top-left cell and move the active selection to K14
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K14").Select()
